$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 327
$ws.Range("C2").Value = 252
$ws.Range("L2").Value = 38
$ws.Range("M2").Value = 28
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 7

$ws.Range("B3").Value = 973
$ws.Range("C3").Value = 914
$ws.Range("D3").Value = 11
$ws.Range("J3").Value = 344
$ws.Range("L3").Value = 108
$ws.Range("M3").Value = 97
$ws.Range("Q3").Value = 27

$ws.Range("B4").Value = 637
$ws.Range("C4").Value = 561
$ws.Range("J4").Value = 16
$ws.Range("L4").Value = 124
$ws.Range("M4").Value = 113
$ws.Range("Q4").Value = 11
$ws.Range("R4").Value = 10

$ws.Range("B5").Value = 310
$ws.Range("C5").Value = 246
$ws.Range("L5").Value = 53
$ws.Range("M5").Value = 44
$ws.Range("Q5").Value = 16

$ws.Range("B6").Value = 150
$ws.Range("C6").Value = 66
$ws.Range("L6").Value = 28
$ws.Range("M6").Value = 14

$ws.Range("B7").Value = 951
$ws.Range("C7").Value = 885
$ws.Range("D7").Value = 14
$ws.Range("J7").Value = 355
$ws.Range("L7").Value = 237
$ws.Range("M7").Value = 231
$ws.Range("Q7").Value = 17
$ws.Range("R7").Value = 15

$ws.Range("B8").Value = 934
$ws.Range("C8").Value = 851
$ws.Range("D8").Value = 2
$ws.Range("J8").Value = 200
$ws.Range("L8").Value = 125
$ws.Range("M8").Value = 115

$ws.Range("B9").Value = 388
$ws.Range("C9").Value = 316
$ws.Range("L9").Value = 86
$ws.Range("M9").Value = 78

$ws.Range("B10").Value = 493
$ws.Range("C10").Value = 437
$ws.Range("J10").Value = 15
$ws.Range("L10").Value = 79
$ws.Range("M10").Value = 75
$ws.Range("Q10").Value = 2

$ws.Range("B11").Value = 453
$ws.Range("C11").Value = 377
$ws.Range("J11").Value = 1
$ws.Range("L11").Value = 105
$ws.Range("M11").Value = 83
$ws.Range("Q11").Value = 22
$ws.Range("R11").Value = 21

$ws.Range("B12").Value = 397
$ws.Range("C12").Value = 299
$ws.Range("L12").Value = 69
$ws.Range("M12").Value = 59
$ws.Range("Q12").Value = 5
$ws.Range("R12").Value = 4

$ws.Range("B13").Value = 110
$ws.Range("C13").Value = 54
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 23
$ws.Range("M13").Value = 13
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 0

$ws.Range("B14").Value = 154
$ws.Range("C14").Value = 87
$ws.Range("L14").Value = 15
$ws.Range("M14").Value = 8

$ws.Range("B15").Value = 736
$ws.Range("C15").Value = 652
$ws.Range("D15").Value = 1
$ws.Range("J15").Value = 56
$ws.Range("L15").Value = 133
$ws.Range("M15").Value = 117
$ws.Range("Q15").Value = 30

$ws.Range("B16").Value = 883
$ws.Range("C16").Value = 832
$ws.Range("D16").Value = 21
$ws.Range("J16").Value = 372
$ws.Range("L16").Value = 110
$ws.Range("M16").Value = 107
$ws.Range("Q16").Value = 24
$ws.Range("R16").Value = 22

$ws.Range("B17").Value = 591
$ws.Range("C17").Value = 533
$ws.Range("J17").Value = 26
$ws.Range("L17").Value = 116
$ws.Range("M17").Value = 114

$ws.Range("B18").Value = 745
$ws.Range("C18").Value = 656
$ws.Range("J18").Value = 113
$ws.Range("L18").Value = 133
$ws.Range("M18").Value = 108

$ws.Range("B19").Value = 623
$ws.Range("C19").Value = 574
$ws.Range("D19").Value = 1
$ws.Range("J19").Value = 39
$ws.Range("L19").Value = 103
$ws.Range("M19").Value = 97

$ws.Range("B20").Value = 508
$ws.Range("C20").Value = 455
$ws.Range("D20").Value = 1
$ws.Range("J20").Value = 17
$ws.Range("L20").Value = 77
$ws.Range("M20").Value = 74
$ws.Range("Q20").Value = 8

$ws.Range("B21").Value = 1142
$ws.Range("C21").Value = 1089
$ws.Range("D21").Value = 16
$ws.Range("J21").Value = 480
$ws.Range("L21").Value = 101
$ws.Range("M21").Value = 103
$ws.Range("Q21").Value = 3
$ws.Range("R21").Value = 1

$ws.Range("B22").Value = 630
$ws.Range("C22").Value = 621
$ws.Range("D22").Value = 9
$ws.Range("J22").Value = 287
$ws.Range("L22").Value = 60
$ws.Range("M22").Value = 56
$ws.Range("Q22").Value = 16

$ws.Range("B23").Value = 429
$ws.Range("C23").Value = 428
$ws.Range("D23").Value = 2
$ws.Range("J23").Value = 93
$ws.Range("L23").Value = 114
$ws.Range("M23").Value = 112

$ws.Range("B24").Value = 388
$ws.Range("C24").Value = 315
$ws.Range("L24").Value = 83
$ws.Range("M24").Value = 73
$ws.Range("Q24").Value = 20
$ws.Range("R24").Value = 19
